# Update "想去人数" (want-to-go count) figures with freshly scraped data.
$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F6").Value  = 2301
$wsExhibit.Range("F16").Value = 435
$wsExhibit.Range("F17").Value = 861
$wsExhibit.Range("F19").Value = 3188
$wsExhibit.Range("F25").Value = 264
$wsExhibit.Range("F30").Value = 831

$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F3").Value = 2942

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F13").Value = 2301
$wsAll.Range("F30").Value = 435
$wsAll.Range("F32").Value = 861
$wsAll.Range("F35").Value = 3188
$wsAll.Range("F40").Value = 264
$wsAll.Range("F50").Value = 831
